$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Populate T_cond [K] (A) and T_evap [K] (B) columns for rows 2-64
$ws.Range("A2").Value = 313
$ws.Range("B2").Value = 257
$ws.Range("A3").Value = 313
$ws.Range("B3").Value = 268.2
$ws.Range("A4").Value = 313
$ws.Range("B4").Value = 283
$ws.Range("A5").Value = 323
$ws.Range("B5").Value = 266.4
$ws.Range("A6").Value = 323
$ws.Range("B6").Value = 278
$ws.Range("A7").Value = 323
$ws.Range("B7").Value = 288
$ws.Range("A8").Value = 326.1
$ws.Range("B8").Value = 257
$ws.Range("A9").Value = 330
$ws.Range("B9").Value = 271.2
$ws.Range("A10").Value = 332
$ws.Range("B10").Value = 283
$ws.Range("A11").Value = 296
$ws.Range("B11").Value = 257
$ws.Range("A12").Value = 313
$ws.Range("B12").Value = 257
$ws.Range("A13").Value = 313
$ws.Range("B13").Value = 268.2
$ws.Range("A14").Value = 323
$ws.Range("B14").Value = 266.4
$ws.Range("A15").Value = 323
$ws.Range("B15").Value = 266.4
$ws.Range("A16").Value = 323.1
$ws.Range("B16").Value = 278
$ws.Range("A17").Value = 326
$ws.Range("B17").Value = 257
$ws.Range("A18").Value = 328
$ws.Range("B18").Value = 251.8
$ws.Range("A19").Value = 333
$ws.Range("B19").Value = 257
$ws.Range("A20").Value = 333.1
$ws.Range("B20").Value = 268.2
$ws.Range("A21").Value = 340
$ws.Range("B21").Value = 263
$ws.Range("A22").Value = 340
$ws.Range("B22").Value = 273.2
$ws.Range("A23").Value = 340
$ws.Range("B23").Value = 285.1
$ws.Range("A24").Value = 296
$ws.Range("B24").Value = 248.2
$ws.Range("A25").Value = 296
$ws.Range("B25").Value = 257
$ws.Range("A26").Value = 313
$ws.Range("B26").Value = 248.2
$ws.Range("A27").Value = 313
$ws.Range("B27").Value = 257
$ws.Range("A28").Value = 313
$ws.Range("B28").Value = 268.2
$ws.Range("A29").Value = 323
$ws.Range("B29").Value = 248.2
$ws.Range("A30").Value = 323
$ws.Range("B30").Value = 266.4
$ws.Range("A31").Value = 323
$ws.Range("B31").Value = 266.4
$ws.Range("A32").Value = 323.1
$ws.Range("B32").Value = 278
$ws.Range("A33").Value = 326
$ws.Range("B33").Value = 257
$ws.Range("A34").Value = 331.2
$ws.Range("B34").Value = 255.1
$ws.Range("A35").Value = 333
$ws.Range("B35").Value = 257
$ws.Range("A36").Value = 333
$ws.Range("B36").Value = 268.2
$ws.Range("A37").Value = 333
$ws.Range("B37").Value = 288
$ws.Range("A38").Value = 340
$ws.Range("B38").Value = 263
$ws.Range("A39").Value = 340
$ws.Range("B39").Value = 273.2
$ws.Range("A40").Value = 340
$ws.Range("B40").Value = 285.1
$ws.Range("A41").Value = 295.8
$ws.Range("B41").Value = 272
$ws.Range("A42").Value = 296
$ws.Range("B42").Value = 257
$ws.Range("A43").Value = 297.2
$ws.Range("B43").Value = 248.2
$ws.Range("A44").Value = 313
$ws.Range("B44").Value = 248.2
$ws.Range("A45").Value = 313
$ws.Range("B45").Value = 257
$ws.Range("A46").Value = 313
$ws.Range("B46").Value = 268.2
$ws.Range("A47").Value = 323
$ws.Range("B47").Value = 248.2
$ws.Range("A48").Value = 323
$ws.Range("B48").Value = 266.4
$ws.Range("A49").Value = 326
$ws.Range("B49").Value = 257
$ws.Range("A50").Value = 327.7
$ws.Range("B50").Value = 278
$ws.Range("A51").Value = 333
$ws.Range("B51").Value = 268.2
$ws.Range("A52").Value = 336
$ws.Range("B52").Value = 259.3
$ws.Range("A53").Value = 336
$ws.Range("B53").Value = 273.2
$ws.Range("A54").Value = 340
$ws.Range("B54").Value = 263
$ws.Range("A55").Value = 340
$ws.Range("B55").Value = 278
$ws.Range("A56").Value = 296.2
$ws.Range("B56").Value = 248.2
$ws.Range("A57").Value = 313
$ws.Range("B57").Value = 248.2
$ws.Range("A58").Value = 313
$ws.Range("B58").Value = 257
$ws.Range("A59").Value = 313
$ws.Range("B59").Value = 268.2
$ws.Range("A60").Value = 323
$ws.Range("B60").Value = 248.2
$ws.Range("A61").Value = 323
$ws.Range("B61").Value = 266.4
$ws.Range("A62").Value = 326
$ws.Range("B62").Value = 257
$ws.Range("A63").Value = 333
$ws.Range("B63").Value = 257
$ws.Range("A64").Value = 333.1
$ws.Range("B64").Value = 268.2

# Update the view/selection state as recorded in the workbook
$ws.Range("C67").Select()
